$wb = $excel.ActiveWorkbook
$ws17 = $wb.Worksheets.Item(16)
$ws17.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollWorkbookTabs(1, 15)
Write-Output "done"
